# Apply updated dSF (column F) values as per repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = -3
    8  = -2
    11 = -4
    13 = -3
    14 = -2
    16 = -11
    17 = -4
    19 = 5
    21 = -7
    22 = -2
    23 = -4
    24 = 1
    25 = -3
    28 = -6
    30 = -1
    31 = -4
    32 = -6
    33 = -2
    34 = -4
    35 = -5
    40 = 1
    41 = 7
    42 = -13
    43 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
